$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "NA" id log entries for maianthemum (row 3) curve, matching the
# pattern already present in row 2 (trillium)
$ws.Range("C3").Value = "NA"
$ws.Range("E3").Value = "NA"

# Update the active selection to reflect where the user was last working
$ws.Range("D3").Select()
